$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '69.655.22'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +4.59%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.608.89'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +4.49%  '
$ws.Range('E4').Value = '  +0.13%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '630.99'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +5.04%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '159.01'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +8.23%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '3.606.57'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +4.54%  '
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('E9').Value = '  +3.94%  '
$ws.Range('E10').Value = '  +10.18%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '7.49'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +8.94%  '
$ws.Range('E12').Value = '  +6.07%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000226'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +6.06%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '33.75'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +9.24%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.222.15'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +4.83%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.610.76'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +5.12%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '69.584.56'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +4.65%  '
$ws.Range('E18').Value = '  +1.22%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.76'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +6.21%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '16.15'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +8.96%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '10.23'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +14.47%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '463.82'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +5.83%  '
$ws.Range('E23').Value = '  +4.62%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '78.95'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +3.06%  '
$ws.Range('E25').Value = '  +10.06%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '10.78'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +8.15%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '3.756.71'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +4.85%  '
$ws.Range('E28').Value = '  -0.06%  '
$ws.Range('E29').Value = '  +14.24%  '
$ws.Range('E30').Value = '  +6.42%  '
$ws.Range('E31').Value = '  +13.45%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.174'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +9.37%  '
$ws.Range('E33').Value = '  +8.63%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.998'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.12%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.97'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +7.25%  '
$ws.Range('E36').Value = '  +4.69%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.608.94'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +5.12%  '
$ws.Range('E38').Value = '  +8.03%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.45'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +15.43%  '
$ws.Range('E40').Value = '  -0.01%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0928'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +8.24%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '179.15'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +3.60%  '
$ws.Range('E43').Value = '  +0.16%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.70'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +6.06%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '31.90'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +23.57%  '
$ws.Range('E46').Value = '  +4.62%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.39'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +14.57%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.78'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +13.08%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '46.06'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.77%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.83'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +4.45%  '
$ws.Range('E51').Value = '  +11.27%  '
